$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in B3: "tactil" -> "tactile"
$ws.Range("B3").Value = "Traiter les informations provenant de l'écran tactile et des boutons"

# Add new function text in B4
$ws.Range("B4").Value = "S'intégrer au boitier en n'altérant pas le design de l'objet"

# Move the active selection to B4
$ws.Range("B4").Select()
